$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.393.26'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.904.67'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +10.67%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '246.63'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '40.68'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.351'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.42%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '52.39'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0725'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.14%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0986'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').Value = '2.182.10'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '12.63'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.36%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.713'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.917.83'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.89'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = '35.391.01'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.98'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.54%  '
$ws.Range('D20').Value = '0.0₃0825'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '241.63'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.63%  '
$ws.Range('E22').Value = '  +2.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.11'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.19%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.31'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '168.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.67'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.90'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.34%  '
$ws.Range('E30').Value = '  +4.68%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.24'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.22%  '
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.19'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.88'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.16%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.917'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.67%  '
$ws.Range('E38').Value = '  +8.28%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '97.61'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.54%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.10'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '16.61'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.68%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0653'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.93%  '
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('D45').Value = '1.359.27'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('B47').Value = 'MultiversX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '46.36'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -6.42%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.42'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.79'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '12.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -5.02%  '
$ws.Range('E51').Value = '  -1.24%  '
